$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 701-702, pushing the existing rows 701-741 down to 703-743.
$ws.Range("A701:R702").EntireRow.Insert()

# Fill the first new row (701) with its data.
$ws.Cells.Item(701, 1).Value = 3
$ws.Cells.Item(701, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(701, 3).Value = "Coquimbo"
$ws.Cells.Item(701, 4).Value = 45041
$ws.Cells.Item(701, 5).Value = 5
$ws.Cells.Item(701, 6).Value = 100112032
$ws.Cells.Item(701, 7).Value = "Zapallo italiano"
$ws.Cells.Item(701, 8).Value = "Sin especificar"
$ws.Cells.Item(701, 9).Value = "Primera"
$ws.Cells.Item(701, 10).Value = 195
$ws.Cells.Item(701, 11).Value = 4000
$ws.Cells.Item(701, 12).Value = 4500
$ws.Cells.Item(701, 13).Value = 4269
$ws.Cells.Item(701, 14).Value = "$/caja 36 unidades"
$ws.Cells.Item(701, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(701, 16).Value = 119
$ws.Cells.Item(701, 17).Value = 36
$ws.Cells.Item(701, 18).Value = "Hortaliza"

# Fill the second new row (702) with its data.
$ws.Cells.Item(702, 1).Value = 3
$ws.Cells.Item(702, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(702, 3).Value = "Coquimbo"
$ws.Cells.Item(702, 4).Value = 45041
$ws.Cells.Item(702, 5).Value = 5
$ws.Cells.Item(702, 6).Value = 100112032
$ws.Cells.Item(702, 7).Value = "Zapallo italiano"
$ws.Cells.Item(702, 8).Value = "Sin especificar"
$ws.Cells.Item(702, 9).Value = "Primera"
$ws.Cells.Item(702, 10).Value = 200
$ws.Cells.Item(702, 11).Value = 7000
$ws.Cells.Item(702, 12).Value = 7500
$ws.Cells.Item(702, 13).Value = 7272
$ws.Cells.Item(702, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(702, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(702, 16).Value = 121
$ws.Cells.Item(702, 17).Value = 60
$ws.Cells.Item(702, 18).Value = "Hortaliza"
